# Fix Training Data Issue (#48)
# The "Date" column (BF) held the sheet's source filename-derived string
# "5-24-2007-08" instead of an actual ISO date. Correct it to "2008-05-24"
# for every data row (rows 2-31; row 1 is the "Date" header).
#
# The literal value must remain a text string, not get auto-converted to a
# serial date by Excel's General-format type inference, so it is entered
# with a leading apostrophe (quote-prefix) exactly like a user forcing text
# entry in the Excel UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Range("BF$r")
    $cell.Value = "'2008-05-24"
}
